$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header label in B1 to the new, more specific prompt text.
$ws.Range("B1").Value = "Your name (including last initial):"

# Update the active selection to reflect where the cursor ended up after editing.
$ws.Range("B1").Select()
